$wb = $excel.ActiveWorkbook

# Sheet R1
$ws = $wb.Worksheets.Item("R1")
$ws.Range("G2").Value = "3922:26:39"
$ws.Range("G3").Value = "61:59:17"
$ws.Range("D5").Value = "LTH2121"
$ws.Range("J5").Value = "Dead"

# Sheet R2
$ws = $wb.Worksheets.Item("R2")
$ws.Range("G2").Value = "12103:50:20"
$ws.Range("G3").Value = "3233:33:49"
$ws.Range("G4").Value = "471:45:23"

# Sheet R4
$ws = $wb.Worksheets.Item("R4")
$ws.Range("G2").Value = "2949:40:09"
$ws.Range("G3").Value = "176:52:24"

# Sheet R5
$ws = $wb.Worksheets.Item("R5")
$ws.Range("G2").Value = "423:39:08"

# Sheet R6
$ws = $wb.Worksheets.Item("R6")
$ws.Range("G2").Value = "64:11:26"
